$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = 1.95
$ws.Range("I9").Value = 4.33
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 7.5
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("U9").Value = 1.95
$ws.Range("V9").Value = 1.8
$ws.Range("W9").Value = 6.5
$ws.Range("X9").Value = 8.5
$ws.Range("Z9").Value = 17
$ws.Range("AH9").Value = 10
$ws.Range("AK9").Value = 41
$ws.Range("AU9").Value = 8.5
